$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: title & link update
$ws.Range("D5").Value = "베르누이 방정식"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/10/Bernoulli_equation.html"

# Row 8: title update
$ws.Range("D8").Value = "카카오브레인"

# Row 9: title & link update
$ws.Range("D9").Value = "수학, 통계학 중심의 AI대학원 커리큘럼"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/curriculum-external-confirmation/#utm_source=rss&utm_medium=rss&utm_campaign=curriculum-external-confirmation"

# Row 29: title update
$ws.Range("D29").Value = "프로메디우스"

# Row 43: title & link update
$ws.Range("D43").Value = "[원격용] 윈도우키, 한영키, 알트탭 매핑 오토핫키"
$ws.Range("E43").Value = "https://nittaku.tistory.com/513"
